# "执行脚本" (Execute Script) was finished and removed from the task list on
# the "事件功能" (Event Function) sheet. Remove its row, which shifts every
# row below it up by one, and mark the seven items that follow it (now rows
# 4-10, previously rows 5-11) as completed ("已完成") to match the other
# finished rows above.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

# Remove the row for "执行脚本" - everything below shifts up one row.
$ws.Rows.Item(4).Delete()

# The next 7 rows (now rows 4-10) are now marked as completed, matching the
# green "已完成" styling used by rows 1-3.
$ws.Range("B4:B10").Value = "已完成"
$ws.Range("B4:B10").Interior.Color = 5287936

# Move the active selection to B10, matching where the edit left off.
$ws.Activate()
$ws.Range("B10").Select()
